# Regenerate the localization-status report: the Overview/zh-cn/de-de "Status"
# cells that previously said "Ready for handoff" now read "In Translation",
# and the now-narrower "Status" columns are resized to match.

$wb = $excel.ActiveWorkbook

# --- 1. Update status text from "Ready for handoff" to "In Translation" ---
# NOTE: compare with the literal on the LEFT side of -eq; some cells hold a
# Boolean value whose .Text reads as "True"/"False", and PowerShell's -eq
# coerces the right-hand side to the left operand's type when the left
# operand is typed first, which can create false positives if reversed.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ("Ready for handoff" -eq [string]$cell.Text) {
            $cell.Value = "In Translation"
        }
    }
}

# --- 2. Narrow the "Status" columns ---
# Overview sheet: columns E (zh-cn) and F (de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn sheet: column C (Status)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# de-de sheet: column C (Status)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
